# Update "想去人数" (want-to-go count) figures for the two events that
# appear on both the "展览" sheet and the "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1956
$wsExhibition.Range("F4").Value = 841
$wsExhibition.Range("F5").Value = 984
$wsExhibition.Range("F6").Value = 344

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1956
$wsAll.Range("F5").Value = 841
$wsAll.Range("F6").Value = 984
$wsAll.Range("F7").Value = 344
